# Fixed naive component forecaster bug - Presentation state 11.02.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear cells that should no longer contain values (removed from the sheet)
$ws.Range("C2").ClearContents()
$ws.Range("E2").ClearContents()
$ws.Range("C3").ClearContents()

# Updated (recalculated) values - tiny floating point precision corrections
$ws.Range("E3").Value = -0.180093374131185
$ws.Range("C4").Value = 0.348613976222456
$ws.Range("C5").Value = -0.1384957661262676
$ws.Range("C6").Value = 1.566479473280191
$ws.Range("C7").Value = 0.7307568962937161
$ws.Range("C8").Value = 0.8188188121642126
$ws.Range("E8").Value = 0.960760217268164
$ws.Range("E9").Value = 1.375398114243231
$ws.Range("C10").Value = 1.9846842782967
$ws.Range("E10").Value = 1.47327408793585
$ws.Range("E11").Value = 1.681032827388385
$ws.Range("C13").Value = 1.064321453542272
$ws.Range("E13").Value = 0.7767182380207682
$ws.Range("C14").Value = 1.361817904277718
$ws.Range("C15").Value = -4.352425014431327
$ws.Range("E15").Value = 0.9348518890383906
$ws.Range("E16").Value = 5.161235657134755
$ws.Range("E17").Value = 2.430255857698516
$ws.Range("C18").Value = -0.9008525709169657
$ws.Range("E18").Value = 1.982587461121321
$ws.Range("C19").Value = 0.2738544794132602
